$d = $word.ActiveDocument

# 1. Update the date field text from 2/8/2025 to 3/25/2025
$d.Content.Find.Execute("2/8/2025", $true, $false, $false, $false, $false, $true, 1, $false, "3/25/2025", 2)

# 2. Reword the EPR / thrust limit sentence.
$old = "EPR is an indicator of the thrust produced by the engine. As a result, it can be used to ensure the engine is operating within safe parameters, which is called the thrust limit."
$new = "EPR is an indicator of the thrust produced by the engine and is used to ensure the engine is operating within safe parameters. This is called the thrust limit."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
